$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 13.59876144578819
$ws.Range("D2").Value = 10.33025063612868
$ws.Range("E2").Value = 15.6201235816481
$ws.Range("F2").Value = 31.20130980610337
$ws.Range("G2").Value = 30.74189977499416
$ws.Range("H2").Value = 14.91809306450685
$ws.Range("I2").Value = 26.56936379662695
$ws.Range("J2").Value = 11.15413369404624
$ws.Range("K2").Value = 9.389141871653191
$ws.Range("L2").Value = 9.347667937842663
$ws.Range("O2").Value = 22.94644218535537
$ws.Range("B3").Value = 13.42491336756516
$ws.Range("D3").Value = 10.33680074722178
$ws.Range("E3").Value = 15.66576440625294
$ws.Range("F3").Value = 31.3046774702009
$ws.Range("G3").Value = 30.86094587252274
$ws.Range("H3").Value = 14.9682581153946
$ws.Range("I3").Value = 26.69695683065929
$ws.Range("J3").Value = 11.17930656722714
$ws.Range("K3").Value = 9.038988737968175
$ws.Range("L3").Value = 9.293022183203908
$ws.Range("O3").Value = 23.03416769223677
$ws.Range("B4").Value = 13.31860683793639
$ws.Range("D4").Value = 10.34217889220578
$ws.Range("E4").Value = 15.6956174040263
$ws.Range("F4").Value = 31.37472359664035
$ws.Range("G4").Value = 30.94283025612953
$ws.Range("H4").Value = 15.00121564842478
$ws.Range("I4").Value = 26.77956808317264
$ws.Range("J4").Value = 11.19558984799751
$ws.Range("K4").Value = 8.815487649723725
$ws.Range("L4").Value = 9.260232783906149
$ws.Range("O4").Value = 23.09245471724982
$ws.Range("B5").Value = 13.27544086479237
$ws.Range("D5").Value = 10.34471239689204
$ws.Range("E5").Value = 15.70824359950783
$ws.Range("F5").Value = 31.40491985603799
$ws.Range("G5").Value = 30.97840148370277
$ws.Range("H5").Value = 15.01518867998672
$ws.Range("I5").Value = 26.81430845814576
$ws.Range("J5").Value = 11.20243397160338
$ws.Range("K5").Value = 8.722361352655946
$ws.Range("L5").Value = 9.247071994819461
$ws.Range("O5").Value = 23.11731859817586
$ws.Range("B6").Value = 13.2682838490571
$ws.Range("D6").Value = 10.34515375340421
$ws.Range("E6").Value = 15.71036802938387
$ws.Range("F6").Value = 31.41003362598413
$ws.Range("G6").Value = 30.98444086851655
$ws.Range("H6").Value = 15.01754167382994
$ws.Range("I6").Value = 26.8201421019381
$ws.Range("J6").Value = 11.20358304665675
$ws.Range("K6").Value = 8.706776997734005
$ws.Range("L6").Value = 9.244899069108683
$ws.Range("O6").Value = 23.12151433153906
$ws.Range("B7").Value = 13.31802400043209
$ws.Range("D7").Value = 10.34221167471339
$ws.Range("E7").Value = 15.69578581817915
$ws.Range("F7").Value = 31.37512414847205
$ws.Range("G7").Value = 30.94330107379307
$ws.Range("H7").Value = 15.0014018961981
$ws.Range("I7").Value = 26.78003224588362
$ws.Range("J7").Value = 11.19568130502535
$ws.Range("K7").Value = 8.814239878975922
$ws.Range("L7").Value = 9.260054466453896
$ws.Range("O7").Value = 23.09278554111435
$ws.Range("B8").Value = 13.53875491477913
$ws.Range("D8").Value = 10.33222810393497
$ws.Range("E8").Value = 15.63548137017385
$ws.Range("F8").Value = 31.23558444186346
$ws.Range("G8").Value = 30.78111815291954
$ws.Range("H8").Value = 14.93494275626057
$ws.Range("I8").Value = 26.61247327477432
$ws.Range("J8").Value = 11.16264200543635
$ws.Range("K8").Value = 9.270225670694177
$ws.Range("L8").Value = 9.328673262438212
$ws.Range("O8").Value = 22.97577142331287
$ws.Range("B9").Value = 13.97310535419625
$ws.Range("D9").Value = 10.32337519581865
$ws.Range("E9").Value = 15.5317010624111
$ws.Range("F9").Value = 31.01423426777663
$ws.Range("G9").Value = 30.53314254704701
$ws.Range("H9").Value = 14.82170402188708
$ws.Range("I9").Value = 26.31765726487633
$ws.Range("J9").Value = 11.10438719451152
$ws.Range("K9").Value = 10.09346281213534
$ws.Range("L9").Value = 9.468875690557152
$ws.Range("O9").Value = 22.78143877624916
$ws.Range("B10").Value = 14.29052569285846
$ws.Range("D10").Value = 10.32335622014657
$ws.Range("E10").Value = 15.46422558501945
$ws.Range("F10").Value = 30.88360513692709
$ws.Range("G10").Value = 30.3940882706026
$ws.Range("H10").Value = 14.7488976319798
$ws.Range("I10").Value = 26.12149608906765
$ws.Range("J10").Value = 11.06553349091945
$ws.Range("K10").Value = 10.65109364631567
$ws.Range("L10").Value = 9.574751740696419
$ws.Range("O10").Value = 22.66012747430565
$ws.Range("B11").Value = 14.43401994283967
$ws.Range("D11").Value = 10.3247430415326
$ws.Range("E11").Value = 15.43542273163673
$ws.Range("F11").Value = 30.83114962525852
$ws.Range("G11").Value = 30.34027673023848
$ws.Range("H11").Value = 14.71802672191015
$ws.Range("I11").Value = 26.03666520703997
$ws.Range("J11").Value = 11.04870673853024
$ws.Range("K11").Value = 10.8938648416299
$ws.Range("L11").Value = 9.623415621448034
$ws.Range("O11").Value = 22.60961085610521
$ws.Range("B12").Value = 14.48818399460891
$ws.Range("D12").Value = 10.32546763204921
$ws.Range("E12").Value = 15.42478707646827
$ws.Range("F12").Value = 30.81228981734013
$ws.Range("G12").Value = 30.32126391070142
$ws.Range("H12").Value = 14.70665975698107
$ws.Range("I12").Value = 26.00517302488946
$ws.Range("J12").Value = 11.04245622049322
$ws.Range("K12").Value = 10.98418368503997
$ws.Range("L12").Value = 9.641904558879606
$ws.Range("O12").Value = 22.59115374815196
$ws.Range("B13").Value = 14.47652723351263
$ws.Range("D13").Value = 10.32530272925992
$ws.Range("E13").Value = 15.4270655976377
$ws.Range("F13").Value = 30.81630693455844
$ws.Range("G13").Value = 30.32529788298864
$ws.Range("H13").Value = 14.70909346646032
$ws.Range("I13").Value = 26.01192737544896
$ws.Range("J13").Value = 11.04379699016579
$ws.Range("K13").Value = 10.96480432607378
$ws.Range("L13").Value = 9.637920098580704
$ws.Range("O13").Value = 22.59509889212974
$ws.Range("B14").Value = 14.43847979105706
$ws.Range("D14").Value = 10.32479866448196
$ws.Range("E14").Value = 15.43454229550568
$ws.Range("F14").Value = 30.82957788638002
$ws.Range("G14").Value = 30.33868515448852
$ws.Range("H14").Value = 14.71708507920068
$ws.Range("I14").Value = 26.03406168178367
$ws.Range("J14").Value = 11.04819007448641
$ws.Range("K14").Value = 10.90132800233738
$ws.Range("L14").Value = 9.624935558361548
$ws.Range("O14").Value = 22.60807889572571
$ws.Range("B15").Value = 14.41515066429854
$ws.Range("D15").Value = 10.32451584411136
$ws.Range("E15").Value = 15.43915730677761
$ws.Range("F15").Value = 30.83783752617252
$ws.Range("G15").Value = 30.34706311333594
$ws.Range("H15").Value = 14.72202225331798
$ws.Range("I15").Value = 26.04770176087444
$ws.Range("J15").Value = 11.05089676065702
$ws.Range("K15").Value = 10.86223555769154
$ws.Range("L15").Value = 9.616989767469846
$ws.Range("O15").Value = 22.61611712690107
$ws.Range("B16").Value = 14.28112603177973
$ws.Range("D16").Value = 10.3232935646493
$ws.Range("E16").Value = 15.46614593110069
$ws.Range("F16").Value = 30.88717361662819
$ws.Range("G16").Value = 30.39779558496677
$ws.Range("H16").Value = 14.75096035538642
$ws.Range("I16").Value = 26.12712844533002
$ws.Range("J16").Value = 11.06665017691505
$ws.Range("K16").Value = 10.6350045786049
$ws.Range("L16").Value = 9.571580574142157
$ws.Range("O16").Value = 22.66352289192037
$ws.Range("B17").Value = 14.19864326456798
$ws.Range("D17").Value = 10.3229002933965
$ws.Range("E17").Value = 15.48318667109942
$ws.Range("F17").Value = 30.91922588764102
$ws.Range("G17").Value = 30.43134200587903
$ws.Range("H17").Value = 14.76928881607303
$ws.Range("I17").Value = 26.17698075175614
$ws.Range("J17").Value = 11.0765311983679
$ws.Range("K17").Value = 10.49277946890062
$ws.Range("L17").Value = 9.543843877761946
$ws.Range("O17").Value = 22.69380130650311
$ws.Range("B18").Value = 14.15111905289002
$ws.Range("D18").Value = 10.32280551982095
$ws.Range("E18").Value = 15.49316618798091
$ws.Range("F18").Value = 30.93831729498637
$ws.Range("G18").Value = 30.45152567796902
$ws.Range("H18").Value = 14.7800425696148
$ws.Range("I18").Value = 26.20606909542346
$ws.Range("J18").Value = 11.08229434913267
$ws.Range("K18").Value = 10.40995189829707
$ws.Range("L18").Value = 9.527938219948393
$ws.Range("O18").Value = 22.7116560055434
$ws.Range("B19").Value = 14.13501540985605
$ws.Range("D19").Value = 10.32279603982964
$ws.Range("E19").Value = 15.49657569762551
$ws.Range("F19").Value = 30.94489388855657
$ws.Range("G19").Value = 30.45851195141996
$ws.Range("H19").Value = 14.78371997417033
$ws.Range("I19").Value = 26.21598917587463
$ws.Range("J19").Value = 11.08425938332524
$ws.Range("K19").Value = 10.38173364737808
$ws.Range("L19").Value = 9.522561372451424
$ws.Range("O19").Value = 22.71777672542465
$ws.Range("B20").Value = 14.20743254016361
$ws.Range("D20").Value = 10.32292856441294
$ws.Range("E20").Value = 15.48135422427872
$ws.Range("F20").Value = 30.91574598254709
$ws.Range("G20").Value = 30.42767891235871
$ws.Range("H20").Value = 14.76731581093096
$ws.Range("I20").Value = 26.17163099087295
$ws.Range("J20").Value = 11.07547108661215
$ws.Range("K20").Value = 10.5080258539064
$ws.Range("L20").Value = 9.546791633018438
$ws.Range("O20").Value = 22.6905326377768
$ws.Range("B21").Value = 14.44966031983071
$ws.Range("D21").Value = 10.3249413182484
$ws.Range("E21").Value = 15.43233884914846
$ws.Range("F21").Value = 30.82565262465332
$ws.Range("G21").Value = 30.33471591020465
$ws.Range("H21").Value = 14.71472897967711
$ws.Range("I21").Value = 26.02754317700202
$ws.Range("J21").Value = 11.04689642837471
$ws.Range("K21").Value = 10.92001665862967
$ws.Range("L21").Value = 9.62874786529809
$ws.Range("O21").Value = 22.60424809543463
$ws.Range("B22").Value = 14.60693571425627
$ws.Range("D22").Value = 10.32741857815437
$ws.Range("E22").Value = 15.40188583606427
$ws.Range("F22").Value = 30.77262399406586
$ws.Range("G22").Value = 30.28191478928957
$ws.Range("H22").Value = 14.68224407572124
$ws.Range("I22").Value = 25.93705310646805
$ws.Range("J22").Value = 11.02892865196459
$ws.Range("K22").Value = 11.17985575518111
$ws.Range("L22").Value = 9.682661319243529
$ws.Range("O22").Value = 22.55177615250001
$ws.Range("B23").Value = 14.52310385508856
$ws.Range("D23").Value = 10.32599054078349
$ws.Range("E23").Value = 15.41799472219101
$ws.Range("F23").Value = 30.80039027245735
$ws.Range("G23").Value = 30.3093658618272
$ws.Range("H23").Value = 14.69940960361488
$ws.Range("I23").Value = 25.98501326760432
$ws.Range("J23").Value = 11.03845383853022
$ws.Range("K23").Value = 11.04205033572249
$ws.Range("L23").Value = 9.653858280074768
$ws.Range("O23").Value = 22.57942237275928
$ws.Range("B24").Value = 14.20345922849417
$ws.Range("D24").Value = 10.32291537394945
$ws.Range("E24").Value = 15.48218210476026
$ws.Range("F24").Value = 30.91731717908493
$ws.Range("G24").Value = 30.42933220213464
$ws.Range("H24").Value = 14.76820713222428
$ws.Range("I24").Value = 26.17404828572412
$ws.Range("J24").Value = 11.07595010636864
$ws.Range("K24").Value = 10.50113626279212
$ws.Range("L24").Value = 9.545458825762859
$ws.Range("O24").Value = 22.69200900988727
$ws.Range("B25").Value = 13.85571352012145
$ws.Range("D25").Value = 10.32462712712716
$ws.Range("E25").Value = 15.55823212156523
$ws.Range("F25").Value = 31.06850459027515
$ws.Range("G25").Value = 30.59268068567665
$ws.Range("H25").Value = 14.85051158825082
$ws.Range("I25").Value = 26.39381297301028
$ws.Range("J25").Value = 11.1194509666786
$ws.Range("K25").Value = 9.878793782624827
$ws.Range("L25").Value = 9.430401315860168
$ws.Range("O25").Value = 22.83024485959526
